# Auto-generated: apply cached-value corrections produced by the scheduled
# market-data runner (see commit message) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2472.8235
$ws.Range("J43").Value = 2873.1667
$ws.Range("L43").Value = 2873.1667
$ws.Range("N43").Value = -3011.1667
$ws.Range("H87").Value = 25180
$ws.Range("J87").Value = 25180
$ws.Range("L87").Value = 25180
$ws.Range("N87").Value = -27676
$ws.Range("H90").Value = 25180
$ws.Range("J90").Value = 25180
$ws.Range("L90").Value = 75540
$ws.Range("N90").Value = -88020
$ws.Range("H129").Value = 936.7778
$ws.Range("J129").Value = 1090.1428
$ws.Range("L129").Value = 3270.4284
$ws.Range("N129").Value = -13270.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 40000
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H61").Value = 1710.8055
$ws.Range("I61").Value = 1147.8148
$ws.Range("K61").Value = 1147.8148
$ws.Range("M61").Value = -935.8148000000001
$ws.Range("H74").Value = 3546.7837
$ws.Range("I74").Value = 3483.4333
$ws.Range("K74").Value = 3483.4333
$ws.Range("M74").Value = -2609.4333
$ws.Range("H77").Value = 3546.7837
$ws.Range("I77").Value = 3483.4333
$ws.Range("K77").Value = 17417.1665
$ws.Range("M77").Value = -13049.1665
$ws.Range("H136").Value = 1710.8055
$ws.Range("I136").Value = 1147.8148
$ws.Range("K136").Value = 3443.4444
$ws.Range("M136").Value = -893.4444000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30832
$ws.Range("H86").Value = 1765.6
$ws.Range("I86").Value = 1739.5555
$ws.Range("K86").Value = 1739.5555
$ws.Range("M86").Value = -616.5554999999999
$ws.Range("H89").Value = 1765.6
$ws.Range("I89").Value = 1739.5555
$ws.Range("K89").Value = 8697.7775
$ws.Range("M89").Value = -3081.7775

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2515.4182
$ws.Range("I31").Value = 1021.41174
$ws.Range("J31").Value = 4934.2856
$ws.Range("K31").Value = 1021.41174
$ws.Range("L31").Value = 4934.2856
$ws.Range("M31").Value = -726.41174
$ws.Range("N31").Value = -5524.2856
$ws.Range("H34").Value = 2515.4182
$ws.Range("I34").Value = 1021.41174
$ws.Range("J34").Value = 4934.2856
$ws.Range("K34").Value = 1021.41174
$ws.Range("L34").Value = 4934.2856
$ws.Range("M34").Value = -819.41174
$ws.Range("N34").Value = -5338.2856
$ws.Range("H107").Value = 959.36365
$ws.Range("I107").Value = 694.125
$ws.Range("J107").Value = 1666.6666
$ws.Range("K107").Value = 694.125
$ws.Range("L107").Value = 1666.6666
$ws.Range("M107").Value = 1225.875
$ws.Range("N107").Value = -5506.6666
$ws.Range("H141").Value = 15603.846
$ws.Range("J141").Value = 15603.846
$ws.Range("L141").Value = 15603.846
$ws.Range("N141").Value = -25963.846

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1744.8182
$ws.Range("J5").Value = 5116.6665
$ws.Range("L5").Value = 15349.9995
$ws.Range("N5").Value = -15573.9995
$ws.Range("H23").Value = 180.14285
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = 201.83333
$ws.Range("K23").Value = 150
$ws.Range("L23").Value = 605.49999
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = -1075.49999
$ws.Range("H37").Value = 200143800
$ws.Range("J37").Value = 200143800
$ws.Range("L37").Value = 600431400
$ws.Range("N37").Value = -600431624
$ws.Range("H132").Value = 2444.8936
$ws.Range("I132").Value = 927.7143
$ws.Range("J132").Value = 3088.5454
$ws.Range("K132").Value = 8349.4287
$ws.Range("L132").Value = 27796.9086
$ws.Range("M132").Value = -5819.4287
$ws.Range("N132").Value = -32856.9086
$ws.Range("H135").Value = 1744.8182
$ws.Range("J135").Value = 5116.6665
$ws.Range("L135").Value = 46049.9985
$ws.Range("N135").Value = -51119.9985
$ws.Range("H139").Value = 2317.516
$ws.Range("I139").Value = 1134.8948
$ws.Range("J139").Value = 4190
$ws.Range("K139").Value = 3404.6844
$ws.Range("L139").Value = 12570
$ws.Range("M139").Value = 1735.3156
$ws.Range("N139").Value = -22850
$ws.Range("H140").Value = 60103.332
$ws.Range("I140").Value = 60103.332
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 180309.996
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -175129.996
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1204.8276
$ws.Range("I113").Value = 1283.4286
$ws.Range("J113").Value = 1131.4667
$ws.Range("K113").Value = 1283.4286
$ws.Range("L113").Value = 1131.4667
$ws.Range("M113").Value = 886.5714
$ws.Range("N113").Value = -5471.4667
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 2679.7742
$ws.Range("I122").Value = 2207.7727
$ws.Range("J122").Value = 3833.5557
$ws.Range("K122").Value = 6623.3181
$ws.Range("L122").Value = 11500.6671
$ws.Range("M122").Value = -4173.3181
$ws.Range("N122").Value = -16400.6671
$ws.Range("H133").Value = 50779.75
$ws.Range("J133").Value = 50779.75
$ws.Range("L133").Value = 50779.75
$ws.Range("N133").Value = -60899.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4995.879
$ws.Range("I40").Value = 4438.615
$ws.Range("J40").Value = 7065.7144
$ws.Range("K40").Value = 4438.615
$ws.Range("L40").Value = 7065.7144
$ws.Range("M40").Value = -4302.615
$ws.Range("N40").Value = -7337.7144
$ws.Range("H122").Value = 3784.0715
$ws.Range("I122").Value = 3139.3044
$ws.Range("J122").Value = 6750
$ws.Range("K122").Value = 9417.913199999999
$ws.Range("L122").Value = 20250
$ws.Range("M122").Value = -6967.913199999999
$ws.Range("N122").Value = -25150
$ws.Range("H130").Value = 49811.25
$ws.Range("J130").Value = 49811.25
$ws.Range("L130").Value = 49811.25
$ws.Range("N130").Value = -59851.25
$ws.Range("H140").Value = 67128.5
$ws.Range("J140").Value = 67128.5
$ws.Range("L140").Value = 67128.5
$ws.Range("N140").Value = -77488.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 670.6957
$ws.Range("I107").Value = 604.7646999999999
$ws.Range("J107").Value = 857.5
$ws.Range("K107").Value = 1814.2941
$ws.Range("L107").Value = 2572.5
$ws.Range("M107").Value = 105.7059000000002
$ws.Range("N107").Value = -6412.5
$ws.Range("H108").Value = 39800
$ws.Range("J108").Value = 39800
$ws.Range("L108").Value = 39800
$ws.Range("N108").Value = -47480
$ws.Range("H121").Value = 28890
$ws.Range("J121").Value = 28890
$ws.Range("L121").Value = 28890
$ws.Range("N121").Value = -32384
$ws.Range("H122").Value = 2924.1428
$ws.Range("I122").Value = 2100.8262
$ws.Range("J122").Value = 4502.1665
$ws.Range("K122").Value = 6302.4786
$ws.Range("L122").Value = 13506.4995
$ws.Range("M122").Value = -3852.4786
$ws.Range("N122").Value = -18406.4995
$ws.Range("H132").Value = 6537435
$ws.Range("I132").Value = 592.0909
$ws.Range("K132").Value = 1776.2727
$ws.Range("M132").Value = 753.7273
